# Update the "想去人数" (want-to-go count) figures and one cover image URL
# for the 展览 and 全部类型 sheets, which hold identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value  = 243
    $ws.Range("F4").Value  = 13218
    $ws.Range("F6").Value  = 230
    $ws.Range("F9").Value  = 168
    $ws.Range("F10").Value = 233
    $ws.Range("F11").Value = 474
    $ws.Range("F12").Value = 8
    $ws.Range("F13").Value = 72
    $ws.Range("F17").Value = 431
    $ws.Range("F18").Value = 5582
    $ws.Range("F19").Value = 112
    $ws.Range("F22").Value = 27
    $ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202404/U2EZscfQ1714448575403.jpeg"
    $ws.Range("F23").Value = 41
    $ws.Range("F25").Value = 163
}
